$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.275.93"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "3.395.64"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.43%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "3.402.46"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "3.985.81"
$ws.Range("E13").Value = "  -0.73%  "

$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("E15").Value = "  -2.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.02%  "

$ws.Range("D17").Value = "64.284.75"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").Value = "3.386.12"
$ws.Range("E18").Value = "  -1.51%  "

$ws.Range("E19").Value = "  -0.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("E22").Value = "  -1.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.542"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.52%  "

$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.83%  "

$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("E30").Value = "  +3.64%  "

$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.40%  "

$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.57%  "

$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("D40").Value = "2.874.12"
$ws.Range("E40").Value = "  -4.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("E42").Value = "  -3.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.95"
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.45%  "

$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.60%  "

$ws.Range("E48").Value = "  +3.76%  "

$ws.Range("E49").Value = "  -1.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.92%  "

$ws.Range("E51").Value = "  -1.08%  "
